$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: CATERIN -> LIZZETH, periodo 2211 -> 2302, valores actualizados
$ws.Range("C16").Value = "1128059537"
$ws.Range("D16").Value = "LIZZETH PAOLA OSPINO GONZALEZ"
$ws.Range("E16").Value = "2302"
$ws.Range("F16").Value = 88000
$ws.Range("G16").Value = 2200000

# Row 17: CATERIN -> LIZZETH, periodo queda 2212, valores actualizados
$ws.Range("C17").Value = "1128059537"
$ws.Range("D17").Value = "LIZZETH PAOLA OSPINO GONZALEZ"
$ws.Range("E17").Value = "2212"
$ws.Range("F17").Value = 88000
$ws.Range("G17").Value = 2200000

# Row 18: LIZZETH -> CATERIN, periodo 2212 -> 2304, valores actualizados
$ws.Range("C18").Value = "1044938732"
$ws.Range("D18").Value = "CATERIN ALEZANDRA MAZA PALOMINO"
$ws.Range("E18").Value = "2304"
$ws.Range("F18").Value = 48000
$ws.Range("G18").Value = 1392000

# Row 19: queda CATERIN, periodo 2301 -> 2302, valor de salario actualizado
$ws.Range("C19").Value = "1044938732"
$ws.Range("D19").Value = "CATERIN ALEZANDRA MAZA PALOMINO"
$ws.Range("E19").Value = "2302"
$ws.Range("G19").Value = 1392000

# Row 20: LIZZETH -> CATERIN, periodo queda 2301, valores actualizados
$ws.Range("C20").Value = "1044938732"
$ws.Range("D20").Value = "CATERIN ALEZANDRA MAZA PALOMINO"
$ws.Range("E20").Value = "2301"
$ws.Range("F20").Value = 48000
$ws.Range("G20").Value = 1392000

# Row 21: CATERIN queda, periodo 2302 -> 2212, G actualizado
$ws.Range("C21").Value = "1044938732"
$ws.Range("D21").Value = "CATERIN ALEZANDRA MAZA PALOMINO"
$ws.Range("E21").Value = "2212"
$ws.Range("G21").Value = 1392000

# Row 22: CATERIN queda, periodo 2304 -> 2211, G actualizado
$ws.Range("C22").Value = "1044938732"
$ws.Range("D22").Value = "CATERIN ALEZANDRA MAZA PALOMINO"
$ws.Range("E22").Value = "2211"
$ws.Range("G22").Value = 1392000
